# Auto-generated edit script applying the Zalera_Profits market-data refresh.
# Updates currentAveragePrice / profit columns (H-N) on each crafting-class sheet
# to match the latest scheduled-runner pull, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1636.8823
$ws.Range("I19").Value = 2484.4285
$ws.Range("J19").Value = 1043.6
$ws.Range("K19").Value = 2484.4285
$ws.Range("L19").Value = 1043.6
$ws.Range("M19").Value = -2309.4285
$ws.Range("N19").Value = -1393.6

$ws.Range("H33").Value = 314.8421
$ws.Range("I33").Value = 334.26666
$ws.Range("K33").Value = 334.26666
$ws.Range("M33").Value = -105.26666

$ws.Range("H107").Value = 166671660
$ws.Range("I107").Value = 250001500
$ws.Range("K107").Value = 250001500
$ws.Range("M107").Value = -249999580

$ws.Range("H132").Value = 1036.907
$ws.Range("I132").Value = 1056.8536
$ws.Range("J132").Value = 628
$ws.Range("K132").Value = 3170.5608
$ws.Range("L132").Value = 1884
$ws.Range("M132").Value = -640.5607999999997
$ws.Range("N132").Value = -6944

$ws.Range("H137").Value = 13895467
$ws.Range("I137").Value = 25001518
$ws.Range("K137").Value = 75004554
$ws.Range("M137").Value = -75002004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32326.324
$ws.Range("I32").Value = 34023.543
$ws.Range("K32").Value = 34023.543
$ws.Range("M32").Value = -33736.543

$ws.Range("H64").Value = 5025000.5
$ws.Range("J64").Value = 5025000.5
$ws.Range("L64").Value = 5025000.5
$ws.Range("N64").Value = -5025496.5

$ws.Range("H67").Value = 5025000.5
$ws.Range("J67").Value = 5025000.5
$ws.Range("L67").Value = 5025000.5
$ws.Range("N67").Value = -5026716.5

$ws.Range("H132").Value = 4678.393
$ws.Range("I132").Value = 3496.5715
$ws.Range("J132").Value = 8223.857
$ws.Range("K132").Value = 10489.7145
$ws.Range("L132").Value = 24671.571
$ws.Range("M132").Value = -7959.7145
$ws.Range("N132").Value = -29731.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 32499
$ws.Range("I26").Value = 32499
$ws.Range("K26").Value = 32499
$ws.Range("M26").Value = -32207

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws.Range("H107").Value = 1648.7826
$ws.Range("I107").Value = 1586.7727
$ws.Range("K107").Value = 1586.7727
$ws.Range("M107").Value = 333.2273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 37857.5
$ws.Range("J28").Value = 37857.5
$ws.Range("L28").Value = 37857.5
$ws.Range("N28").Value = -38347.5

$ws.Range("H58").Value = 4706
$ws.Range("I58").Value = 2843.0557
$ws.Range("K58").Value = 2843.0557
$ws.Range("M58").Value = -2640.0557

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 48465.93
$ws.Range("I132").Value = 3927.0908
$ws.Range("K132").Value = 11781.2724
$ws.Range("M132").Value = -9251.2724

$ws.Range("H134").Value = 7006.1665
$ws.Range("I134").Value = 6933.476
$ws.Range("J134").Value = 7515
$ws.Range("K134").Value = 20800.428
$ws.Range("L134").Value = 22545
$ws.Range("M134").Value = -18265.428
$ws.Range("N134").Value = -27615

$ws.Range("H136").Value = 4706
$ws.Range("I136").Value = 2843.0557
$ws.Range("K136").Value = 8529.167099999999
$ws.Range("M136").Value = -5979.167099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1499.8
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H140").Value = 1224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 329.16666
$ws.Range("I22").Value = 215.33333
$ws.Range("J22").Value = 443
$ws.Range("K22").Value = 215.33333
$ws.Range("L22").Value = 443
$ws.Range("M22").Value = 313.66667
$ws.Range("N22").Value = -1501

$ws.Range("H33").Value = 11799.6
$ws.Range("I33").Value = 8090.909
$ws.Range("K33").Value = 8090.909
$ws.Range("M33").Value = -7838.909

$ws.Range("H41").Value = 3357.75
$ws.Range("I41").Value = 999.8570999999999
$ws.Range("K41").Value = 999.8570999999999
$ws.Range("M41").Value = -644.8570999999999

$ws.Range("H132").Value = 8605.294
$ws.Range("I132").Value = 7419.3335
$ws.Range("K132").Value = 22258.0005
$ws.Range("M132").Value = -19728.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1760.4445
$ws.Range("I4").Value = 403
$ws.Range("J4").Value = 2439.1667
$ws.Range("K4").Value = 403
$ws.Range("L4").Value = 2439.1667
$ws.Range("M4").Value = -290
$ws.Range("N4").Value = -2665.1667

$ws.Range("H22").Value = 3829.4644
$ws.Range("I22").Value = 2200.6155
$ws.Range("J22").Value = 5241.1333
$ws.Range("K22").Value = 2200.6155
$ws.Range("L22").Value = 5241.1333
$ws.Range("M22").Value = -1905.6155
$ws.Range("N22").Value = -5831.1333

$ws.Range("H27").Value = 3829.4644
$ws.Range("I27").Value = 2200.6155
$ws.Range("J27").Value = 5241.1333
$ws.Range("K27").Value = 2200.6155
$ws.Range("L27").Value = 5241.1333
$ws.Range("M27").Value = -2093.6155
$ws.Range("N27").Value = -5455.1333

$ws.Range("H28").Value = 1760.4445
$ws.Range("I28").Value = 403
$ws.Range("J28").Value = 2439.1667
$ws.Range("K28").Value = 403
$ws.Range("L28").Value = 2439.1667
$ws.Range("M28").Value = -171
$ws.Range("N28").Value = -2903.1667

$ws.Range("H37").Value = 1760.4445
$ws.Range("I37").Value = 403
$ws.Range("J37").Value = 2439.1667
$ws.Range("K37").Value = 403
$ws.Range("L37").Value = 2439.1667
$ws.Range("M37").Value = -296
$ws.Range("N37").Value = -2653.1667

$ws.Range("H46").Value = 7145.2915
$ws.Range("I46").Value = 1274.5
$ws.Range("J46").Value = 8319.450000000001
$ws.Range("K46").Value = 1274.5
$ws.Range("L46").Value = 8319.450000000001
$ws.Range("M46").Value = -1086.5
$ws.Range("N46").Value = -8695.450000000001

$ws.Range("H55").Value = 712.25
$ws.Range("I55").Value = 886.8
$ws.Range("K55").Value = 886.8
$ws.Range("M55").Value = -713.8

$ws.Range("H68").Value = 2625.8572
$ws.Range("I68").Value = 2736.4
$ws.Range("K68").Value = 2736.4
$ws.Range("M68").Value = -1987.4

$ws.Range("H71").Value = 2625.8572
$ws.Range("I71").Value = 2736.4
$ws.Range("K71").Value = 13682
$ws.Range("M71").Value = -9938

$ws.Range("H132").Value = 10754.55
$ws.Range("I132").Value = 8018.846
$ws.Range("J132").Value = 15835.143
$ws.Range("K132").Value = 24056.538
$ws.Range("L132").Value = 47505.429
$ws.Range("M132").Value = -21526.538
$ws.Range("N132").Value = -52565.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 34632.332
$ws.Range("J63").Value = 34632.332
$ws.Range("L63").Value = 34632.332
$ws.Range("N63").Value = -35880.332

$ws.Range("H66").Value = 34632.332
$ws.Range("J66").Value = 34632.332
$ws.Range("L66").Value = 103896.996
$ws.Range("N66").Value = -110136.996

$ws.Range("H113").Value = 677.9231
$ws.Range("I113").Value = 550.125
$ws.Range("J113").Value = 882.4
$ws.Range("K113").Value = 1650.375
$ws.Range("L113").Value = 2647.2
$ws.Range("M113").Value = 519.625
$ws.Range("N113").Value = -6987.2

$ws.Range("H123").Value = 64690
$ws.Range("J123").Value = 64690
$ws.Range("L123").Value = 64690
$ws.Range("N123").Value = -74490
